$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 19

# --- Capture current (pre-edit) values of columns G, H, I for the data rows ---
# Column layout before the edit:
#   G = MT brut de loyer
#   H = MT brut d'avance
#   I = Taux de taxe
# Column layout after the edit (a "Taux de taxe" column is inserted before the
# old "MT brut de loyer" column, i.e. the misplaced I column moves to G and the
# other two slide right):
#   G = Taux de taxe   (was I)
#   H = MT brut de loyer (was G)
#   I = MT brut d'avance (was H)

$oldG = @{}
$oldH = @{}
$oldI = @{}

for ($r = 2; $r -le $lastRow; $r++) {
    $oldG[$r] = $ws.Cells.Item($r, 7).Value2
    $oldH[$r] = $ws.Cells.Item($r, 8).Value2
    $oldI[$r] = $ws.Cells.Item($r, 9).Value2
}

# --- Header row ---
$ws.Range("G1").Value = "Taux de taxe"
$ws.Range("H1").Value = "MT brut de loyer"
$ws.Range("I1").Value = "MT brut d'avance"

# --- Data rows: shuffle values into their corrected columns ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = $oldI[$r]
    $ws.Cells.Item($r, 8).Value = $oldG[$r]
    $ws.Cells.Item($r, 9).Value = $oldH[$r]
}

# --- New totals row ("Afficher les totaux des montants") ---
# Written as plain literal totals (not live formulas) to mirror the source
# workbook, which stores every amount as a static <v>.
$totalsRow = $lastRow + 1

for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item($totalsRow, $c).Value = " "
}

for ($c = 8; $c -le 13; $c++) {
    $sum = 0.0
    for ($r = 2; $r -le $lastRow; $r++) {
        $cellVal = $ws.Cells.Item($r, $c).Value2
        # SUM() ignores text cells (e.g. the "--" placeholders) just like Excel does.
        if ($cellVal -is [double] -or $cellVal -is [int]) {
            $sum += [double]$cellVal
        }
    }
    $ws.Cells.Item($totalsRow, $c).Value = $sum
}
